# "Generate Report for Handback"
#
# For each language sheet (zh-cn, de-de) the handoff rows (2 & 3) are
# updated to reflect that a handback has now happened:
#   - Status moves from "Ready for handoff" to
#     "Handed back: in sync with en-US" (this text also shows on the
#     Overview sheet, which shares the same cell text).
#   - The "Latest Target File" (E) and "Latest Handback File" (F) columns
#     get populated with links to the same files referenced by the
#     "Source File Name" (A) and "Latest Handoff File" (C) columns.
#   - The "Latest Handback DateTime" (G) column is stamped with the
#     handback time.

$wb = $excel.ActiveWorkbook

function Get-HyperlinkAddress($ws, $cellRef) {
    $target = $ws.Range($cellRef).Address()
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $target) {
            return $hl.Address
        }
    }
    return $null
}

# 1) Update the status text everywhere it appears (Overview + both language
#    sheets) - changing the shared text updates every cell that uses it.
foreach ($sheetName in @("Overview", "zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $used = $ws.UsedRange
    for ($r = 1; $r -le $used.Rows.Count; $r++) {
        for ($c = 1; $c -le $used.Columns.Count; $c++) {
            $cell = $ws.Cells.Item($r, $c)
            if ($cell.Value() -eq "Ready for handoff") {
                $cell.Value = "Handed back: in sync with en-US"
            }
        }
    }
}

# 2) Per language sheet: populate the Target File / Handback File columns
#    and stamp the handback datetime for the two handoff rows.
$handbackInfo = @{
    "zh-cn" = "2016-03-10 09:43:42"
    "de-de" = "2016-03-10 09:44:01"
}

foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $handbackDateTime = $handbackInfo[$sheetName]

    foreach ($row in @(2, 3)) {
        $sourceRef = "A$row"
        $handoffFileRef = "C$row"
        $targetRef = "E$row"
        $handbackRef = "F$row"
        $handbackDateRef = "G$row"

        $sourceDisplay = $ws.Range($sourceRef).Value()
        $sourceAddress = Get-HyperlinkAddress $ws $sourceRef

        $handoffDisplay = $ws.Range($handoffFileRef).Value()
        $handoffAddress = Get-HyperlinkAddress $ws $handoffFileRef

        $ws.Hyperlinks.Add($ws.Range($targetRef), $sourceAddress, "", "", $sourceDisplay) | Out-Null
        $ws.Hyperlinks.Add($ws.Range($handbackRef), $handoffAddress, "", "", $handoffDisplay) | Out-Null

        $ws.Range($handbackDateRef).Value = $handbackDateTime
    }
}
